$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2: refactor the "lowest working prescale" formula so it no longer
#     special-cases F22=0 (which can't happen since F-column defaults to
#     10000 as a sentinel) and instead checks for the real sentinel value.
$ws2.Range("I7").Formula = "=IF(F22=10000,C22,MIN(F7:F22))"

# --- Sheet1: bump the clock frequency input (fclk, MHz) from 10 to 50 ---
$ws1.Range("D9").Value = 50

# --- Sheet1: the old scratch rows used by the timer-module refactor are no
#     longer needed now that the instruction memory moved inside the
#     processor sheet; drop them so the sheet shrinks back down.
$ws1.Rows.Item(19).Resize(2).Delete()

# --- restore selections on both sheets, and make Sheet2 the active tab ---
$ws1.Activate()
$ws1.Range("D10").Select()

$ws2.Activate()
$ws2.Range("I8").Select()
